$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 271, shifting existing rows 271-293 down to 272-294.
$ws.Rows(271).Insert()

# Populate the newly inserted row 271 with the new record.
$ws.Cells.Item(271, 1).Value2 = 7
$ws.Cells.Item(271, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(271, 3).Value2 = "Ñuble"
$ws.Cells.Item(271, 4).Value2 = 45013
$ws.Cells.Item(271, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(271, 5).Value2 = 16
$ws.Cells.Item(271, 6).Value2 = 100112032
$ws.Cells.Item(271, 7).Value2 = "Zapallo italiano"
$ws.Cells.Item(271, 8).Value2 = "Sin especificar"
$ws.Cells.Item(271, 9).Value2 = "Primera"
$ws.Cells.Item(271, 10).Value2 = 60
$ws.Cells.Item(271, 11).Value2 = 6000
$ws.Cells.Item(271, 12).Value2 = 6000
$ws.Cells.Item(271, 13).Value2 = 6000
$ws.Cells.Item(271, 14).Value2 = "`$/caja 50 unidades"
$ws.Cells.Item(271, 15).Value2 = "Región del Maule"
$ws.Cells.Item(271, 16).Value2 = 120
$ws.Cells.Item(271, 17).Value2 = 50
$ws.Cells.Item(271, 18).Value2 = "Hortaliza"
